$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = -0
$ws.Cells.Item(2, 2).Value = -0.0736143064681548
$ws.Cells.Item(2, 3).Value = -0
$ws.Cells.Item(2, 4).Value = 0.2023283625086515
$ws.Cells.Item(2, 5).Value = 0.005361901149070607
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 9).Value = -0
$ws.Cells.Item(2, 10).Value = -0
$ws.Cells.Item(2, 11).Value = 0.005910374655943606
$ws.Cells.Item(2, 12).Value = -0
$ws.Cells.Item(2, 13).Value = 0.2017405004068997
$ws.Cells.Item(2, 14).Value = -0.003175673222564392
$ws.Cells.Item(2, 18).Value = -0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 20).Value = -0.0836551046379089
$ws.Cells.Item(2, 22).Value = 0.01416007321150566
$ws.Cells.Item(2, 23).Value = -0.03298341659304817
$ws.Cells.Item(2, 25).Value = -0
$ws.Cells.Item(2, 26).Value = -0
$ws.Cells.Item(2, 28).Value = 0
$ws.Cells.Item(2, 29).Value = -0.05416417955287071
$ws.Cells.Item(2, 30).Value = 0
$ws.Cells.Item(2, 31).Value = -0.01788162495550331
$ws.Cells.Item(2, 32).Value = 0.0002069792777307436
$ws.Cells.Item(2, 33).Value = -0
$ws.Cells.Item(2, 34).Value = -0
$ws.Cells.Item(2, 35).Value = -0
$ws.Cells.Item(2, 36).Value = 0
$ws.Cells.Item(2, 37).Value = -0
$ws.Cells.Item(2, 38).Value = -0.03184932082569965
$ws.Cells.Item(2, 39).Value = 0
$ws.Cells.Item(2, 40).Value = 0.02827770634814052
$ws.Cells.Item(2, 41).Value = 0.06944358562979185
$ws.Cells.Item(2, 43).Value = 0
$ws.Cells.Item(2, 44).Value = -0
$ws.Cells.Item(2, 46).Value = 0
$ws.Cells.Item(2, 47).Value = -0.1497027310705481
$ws.Cells.Item(2, 49).Value = 0.07080831603100772
$ws.Cells.Item(2, 50).Value = -0.001817919973327277
$ws.Cells.Item(2, 51).Value = -0
$ws.Cells.Item(2, 55).Value = -0
$ws.Cells.Item(2, 56).Value = -0.01418805710578807
$ws.Cells.Item(2, 58).Value = 0.08649163433815991
$ws.Cells.Item(2, 59).Value = 0.03283123518905573
$ws.Cells.Item(2, 62).Value = -0
$ws.Cells.Item(2, 64).Value = 0
$ws.Cells.Item(2, 65).Value = 0.03190481457958391
$ws.Cells.Item(2, 67).Value = -0.04257541708426302
$ws.Cells.Item(2, 68).Value = -0.08905310676590357
$ws.Cells.Item(2, 73).Value = 0
$ws.Cells.Item(2, 74).Value = -0.04640710802875297
$ws.Cells.Item(2, 75).Value = 0
$ws.Cells.Item(2, 76).Value = 0.01094127294829059
$ws.Cells.Item(2, 77).Value = -0.02004983166574863
$ws.Cells.Item(2, 78).Value = -0
$ws.Cells.Item(2, 80).Value = 0
$ws.Cells.Item(2, 82).Value = -0
$ws.Cells.Item(2, 83).Value = 0.03274941356648393
$ws.Cells.Item(2, 85).Value = -0.03099117283538349
$ws.Cells.Item(2, 86).Value = 0.01593166186343111
$ws.Cells.Item(2, 88).Value = -0
$ws.Cells.Item(2, 91).Value = -0
$ws.Cells.Item(2, 92).Value = -0.01061381960660221
$ws.Cells.Item(2, 94).Value = 0.02133540246658532
$ws.Cells.Item(2, 95).Value = 0.03716094318380431
$ws.Cells.Item(2, 98).Value = 0
$ws.Cells.Item(2, 99).Value = -0
$ws.Cells.Item(2, 100).Value = -0
$ws.Cells.Item(2, 101).Value = 0.04597314766486385
$ws.Cells.Item(2, 103).Value = -0.03342648399499332
$ws.Cells.Item(2, 104).Value = 0.01017704690408558
$ws.Cells.Item(2, 108).Value = -0
$ws.Cells.Item(2, 109).Value = -0
$ws.Cells.Item(2, 110).Value = 0.02874569132567836
$ws.Cells.Item(2, 112).Value = 0.02900526664094873
$ws.Cells.Item(2, 113).Value = 0.03380050877759293
$ws.Cells.Item(2, 114).Value = 0
$ws.Cells.Item(2, 115).Value = -0
$ws.Cells.Item(2, 116).Value = -0
$ws.Cells.Item(2, 118).Value = 0
$ws.Cells.Item(2, 119).Value = -0.01950247745448723
$ws.Cells.Item(2, 120).Value = -0
$ws.Cells.Item(2, 121).Value = 0.03637034262361485
$ws.Cells.Item(2, 122).Value = -0.01945341551444906
$ws.Cells.Item(2, 123).Value = -0
$ws.Cells.Item(2, 127).Value = 0
$ws.Cells.Item(2, 128).Value = -0.05668836815106189
$ws.Cells.Item(2, 129).Value = -0
$ws.Cells.Item(2, 130).Value = -0.008248668484950638
$ws.Cells.Item(2, 131).Value = -0.02434073422596091
$ws.Cells.Item(2, 132).Value = 0
$ws.Cells.Item(2, 136).Value = -0
$ws.Cells.Item(2, 137).Value = 0.04085971144248264
$ws.Cells.Item(2, 139).Value = 0.06686095049629477
$ws.Cells.Item(2, 140).Value = -0.02407360759003618
$ws.Cells.Item(2, 145).Value = 0
$ws.Cells.Item(2, 146).Value = 0.04575396385905522
$ws.Cells.Item(2, 147).Value = 0
$ws.Cells.Item(2, 148).Value = -0.0349068518066118
$ws.Cells.Item(2, 149).Value = 0.03654119765287879
$ws.Cells.Item(2, 150).Value = 0
$ws.Cells.Item(2, 151).Value = -0
$ws.Cells.Item(2, 152).Value = 0
$ws.Cells.Item(2, 154).Value = 0
$ws.Cells.Item(2, 155).Value = 0.04340351386436194
$ws.Cells.Item(2, 157).Value = -0.02669773779825179
$ws.Cells.Item(2, 158).Value = 0.01710822948871973
$ws.Cells.Item(2, 160).Value = -0
$ws.Cells.Item(2, 163).Value = -0
$ws.Cells.Item(2, 164).Value = 0.001691558334483329
$ws.Cells.Item(2, 165).Value = 0
$ws.Cells.Item(2, 166).Value = -0.006288705109680439
$ws.Cells.Item(2, 167).Value = -0.008476023910300627
$ws.Cells.Item(2, 168).Value = -0
$ws.Cells.Item(2, 170).Value = -0
$ws.Cells.Item(2, 172).Value = -0
$ws.Cells.Item(2, 173).Value = -0.01406165917757605
$ws.Cells.Item(2, 174).Value = -0
$ws.Cells.Item(2, 175).Value = -0.0184191401413391
$ws.Cells.Item(2, 176).Value = 0.005010644615690384
$ws.Cells.Item(2, 178).Value = -0
$ws.Cells.Item(2, 179).Value = -0
$ws.Cells.Item(2, 181).Value = 0
$ws.Cells.Item(2, 182).Value = -0.03040407898268115
$ws.Cells.Item(2, 184).Value = 0.03399166782388464
$ws.Cells.Item(2, 186).Value = 0
$ws.Cells.Item(2, 187).Value = -0
